# Add support for nullables
# - Adds two new "array" columns (F: ArrayValue1, G: ArrayValue2) with header text
#   and sample values in row 3.
# - Converts the previously-numeric B3/D3/E3 sample cells into a string, a
#   boolean and a real (formatted) date respectively, matching the richer
#   set of supported/nullable CLR types exercised by the test fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two extra columns.
$ws.Range("F1").Value = "ArrayValue1"
$ws.Range("G1").Value = "ArrayValue2"

# Row 3 sample values.
$ws.Range("B3").Value = "string"
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = "Test"

# Give E3 a real date (formatted with a built-in short-date number format,
# numFmtId 14) instead of the plain integer it held before.
$ws.Range("E3").NumberFormat = "mm-dd-yy"
$ws.Range("E3").Value = (Get-Date -Year 2017 -Month 7 -Day 5 -Hour 0 -Minute 0 -Second 0)

$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2

# Move the active selection from F1 to F2, as in the target workbook.
[void]$ws.Range("F2").Select()
